# Daily attendance processing - normalise "Recorded By" (column G) ordering.
# Moves the "System" token that currently leads the list to its
# correctly-ordered position among the other recorder names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "system, System, backup@backdoor.com"
    }
}
